# Generate Report for handoff
# Adds the new source file "f2c27ab1-8f74-4493-98ba-a7649837cbf6.md" as a new
# row just before the trailing ".localization-config" row on every sheet
# (Overview, zh-cn, de-de), shifting that row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" (File Name / zh-cn / de-de summary)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks so none keep stale row references once rows move.
$ws1.Range("A1:C20").Hyperlinks.Delete()

# Push the ".localization-config" row (row 8) down to make room for the new row.
$ws1.Rows.Item(8).Insert()

$ws1.Range("A8").Value = "f2c27ab1-8f74-4493-98ba-a7649837cbf6.md"
$ws1.Range("B8").Value = "Ready for handoff"
$ws1.Range("C8").Value = "Ready for handoff"

$ws1.Range("A9").Value = ".localization-config"
$ws1.Range("B9").Value = "Not to be localized"
$ws1.Range("C9").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/5ec1b24284cb0d06cc9ff951f8a5f627d912ce71/e2e/0ed12709-d088-4d8c-8475-0e19163a68f3.md", "", "", "0ed12709-d088-4d8c-8475-0e19163a68f3.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/e2e/88b63550-690b-4866-9085-6bae40f52bf6.md", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/7fec131d44aeddbeab3ed95d4f3fc61704bc5d08/e2e/c581042a-d6d5-4f26-980c-3c1b59453863.md", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/e2e/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(6,1), "https://github.com/OpenLocalizationTest/oltest/blob/7398222819396b38b5d13c0273fedcf09cf7355b/e2e/45cefd5f-f1e5-46e1-9604-137c12761e97.md", "", "", "45cefd5f-f1e5-46e1-9604-137c12761e97.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(7,1), "https://github.com/OpenLocalizationTest/oltest/blob/fe875b291df0650878992d222bdcbea4ca0cef45/e2e/a5d079c2-1396-4f01-82f7-80769c44f640.md", "", "", "a5d079c2-1396-4f01-82f7-80769c44f640.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(8,1), "https://github.com/OpenLocalizationTest/oltest/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/e2e/f2c27ab1-8f74-4493-98ba-a7649837cbf6.md", "", "", "f2c27ab1-8f74-4493-98ba-a7649837cbf6.md")
$ws1.Hyperlinks.Add($ws1.Cells.Item(9,1), "https://github.com/OpenLocalizationTest/oltest/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn" detail sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A1:I20").Hyperlinks.Delete()

$ws2.Rows.Item(8).Insert()

$ws2.Range("A8").Value = "f2c27ab1-8f74-4493-98ba-a7649837cbf6.md"
$ws2.Range("B8").Value = "Ready for handoff"
$ws2.Range("C8").Value = "f2c27ab1-8f74-4493-98ba-a7649837cbf6.a45b7595bce6941660e6dd85dc4c102f066ae408.zh-cn.xlf"
$ws2.Range("D8").Value = "2016-01-25 03:17:58"
$ws2.Range("G8").Value = "0001-01-01 00:00:00"
$ws2.Range("H8").Value = "Include"

$ws2.Range("A9").Value = ".localization-config"
$ws2.Range("B9").Value = "Not to be localized"
$ws2.Range("D9").Value = "0001-01-01 00:00:00"
$ws2.Range("G9").Value = "0001-01-01 00:00:00"
$ws2.Range("H9").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/5ec1b24284cb0d06cc9ff951f8a5f627d912ce71/e2e/0ed12709-d088-4d8c-8475-0e19163a68f3.md", "", "", "0ed12709-d088-4d8c-8475-0e19163a68f3.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8acabd0f7e67133e87bfe468fadd8ae9ac754cf4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/0ed12709-d088-4d8c-8475-0e19163a68f3.bc16ce64bd8926fc2a9dfb5ce635faaa02561a40.zh-cn.xlf", "", "", "0ed12709-d088-4d8c-8475-0e19163a68f3.bc16ce64bd8926fc2a9dfb5ce635faaa02561a40.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/e2e/88b63550-690b-4866-9085-6bae40f52bf6.md", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0920e71b8e6468777c03d1a93dacdbeffb2766ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.zh-cn.xlf", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/7fec131d44aeddbeab3ed95d4f3fc61704bc5d08/e2e/c581042a-d6d5-4f26-980c-3c1b59453863.md", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50aa5c7e101fd088d78ace95bd7a742aacfe42c4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.zh-cn.xlf", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,5), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6511e9776fae83d47b5c48f8491d88cb96a5fea9/e2e/c581042a-d6d5-4f26-980c-3c1b59453863.md", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,6), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b9a8555989489aabd42a3db540d88defbdee7249/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/qimu/c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.zh-cn.xlf", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/e2e/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0920e71b8e6468777c03d1a93dacdbeffb2766ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.zh-cn.xlf", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(6,1), "https://github.com/OpenLocalizationTest/oltest/blob/7398222819396b38b5d13c0273fedcf09cf7355b/e2e/45cefd5f-f1e5-46e1-9604-137c12761e97.md", "", "", "45cefd5f-f1e5-46e1-9604-137c12761e97.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(6,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f4a6428f9242323a13ee53d8173e7a5e44de5893/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/45cefd5f-f1e5-46e1-9604-137c12761e97.ff005cdba3b917e92e3a7ff3db8ae7ae7f76afa5.zh-cn.xlf", "", "", "45cefd5f-f1e5-46e1-9604-137c12761e97.ff005cdba3b917e92e3a7ff3db8ae7ae7f76afa5.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(7,1), "https://github.com/OpenLocalizationTest/oltest/blob/fe875b291df0650878992d222bdcbea4ca0cef45/e2e/a5d079c2-1396-4f01-82f7-80769c44f640.md", "", "", "a5d079c2-1396-4f01-82f7-80769c44f640.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(7,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34b289b06f96cfd5dad03e058a757c80b9ff9a57/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/a5d079c2-1396-4f01-82f7-80769c44f640.d99f3028deb065a2bc07474fd2c3b3d2b380953e.zh-cn.xlf", "", "", "a5d079c2-1396-4f01-82f7-80769c44f640.d99f3028deb065a2bc07474fd2c3b3d2b380953e.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(8,1), "https://github.com/OpenLocalizationTest/oltest/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/e2e/f2c27ab1-8f74-4493-98ba-a7649837cbf6.md", "", "", "f2c27ab1-8f74-4493-98ba-a7649837cbf6.md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(8,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/f2c27ab1-8f74-4493-98ba-a7649837cbf6.a45b7595bce6941660e6dd85dc4c102f066ae408.zh-cn.xlf", "", "", "f2c27ab1-8f74-4493-98ba-a7649837cbf6.a45b7595bce6941660e6dd85dc4c102f066ae408.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Cells.Item(9,1), "https://github.com/OpenLocalizationTest/oltest/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 3: "de-de" detail sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A1:I20").Hyperlinks.Delete()

$ws3.Rows.Item(8).Insert()

$ws3.Range("A8").Value = "f2c27ab1-8f74-4493-98ba-a7649837cbf6.md"
$ws3.Range("B8").Value = "Ready for handoff"
$ws3.Range("C8").Value = "f2c27ab1-8f74-4493-98ba-a7649837cbf6.a45b7595bce6941660e6dd85dc4c102f066ae408.de-de.xlf"
$ws3.Range("D8").Value = "2016-01-25 03:18:09"
$ws3.Range("G8").Value = "0001-01-01 00:00:00"
$ws3.Range("H8").Value = "Include"

$ws3.Range("A9").Value = ".localization-config"
$ws3.Range("B9").Value = "Not to be localized"
$ws3.Range("D9").Value = "0001-01-01 00:00:00"
$ws3.Range("G9").Value = "0001-01-01 00:00:00"
$ws3.Range("H9").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/5ec1b24284cb0d06cc9ff951f8a5f627d912ce71/e2e/0ed12709-d088-4d8c-8475-0e19163a68f3.md", "", "", "0ed12709-d088-4d8c-8475-0e19163a68f3.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac40d6786f456e05fa6d3bb701e48d0dc27bbaba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/0ed12709-d088-4d8c-8475-0e19163a68f3.bc16ce64bd8926fc2a9dfb5ce635faaa02561a40.de-de.xlf", "", "", "0ed12709-d088-4d8c-8475-0e19163a68f3.bc16ce64bd8926fc2a9dfb5ce635faaa02561a40.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/e2e/88b63550-690b-4866-9085-6bae40f52bf6.md", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d61823019d4254cb5dc26a4b774f1312884e132/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.de-de.xlf", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/7fec131d44aeddbeab3ed95d4f3fc61704bc5d08/e2e/c581042a-d6d5-4f26-980c-3c1b59453863.md", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3952f77917eed0011a9ecea2fa7f2de06a516c6e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.de-de.xlf", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,5), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7516365e0c41fcd2a55de264de05d63d369793a5/e2e/c581042a-d6d5-4f26-980c-3c1b59453863.md", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,6), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f8555fd5caef4a5218d7e7e109b874907a461913/ol-handback/OpenLocalizationTestOrg/oltest.de-de/qimu/c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.de-de.xlf", "", "", "c581042a-d6d5-4f26-980c-3c1b59453863.456f53ebdd5870a037cb78b92a0cb5b7c05fbbe0.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/e2e/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d61823019d4254cb5dc26a4b774f1312884e132/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.de-de.xlf", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(6,1), "https://github.com/OpenLocalizationTest/oltest/blob/7398222819396b38b5d13c0273fedcf09cf7355b/e2e/45cefd5f-f1e5-46e1-9604-137c12761e97.md", "", "", "45cefd5f-f1e5-46e1-9604-137c12761e97.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(6,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/420bc3519b0232c39bb68b41222ab3dfd7d2cf37/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/45cefd5f-f1e5-46e1-9604-137c12761e97.ff005cdba3b917e92e3a7ff3db8ae7ae7f76afa5.de-de.xlf", "", "", "45cefd5f-f1e5-46e1-9604-137c12761e97.ff005cdba3b917e92e3a7ff3db8ae7ae7f76afa5.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(7,1), "https://github.com/OpenLocalizationTest/oltest/blob/fe875b291df0650878992d222bdcbea4ca0cef45/e2e/a5d079c2-1396-4f01-82f7-80769c44f640.md", "", "", "a5d079c2-1396-4f01-82f7-80769c44f640.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(7,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2308e4fec0582e9f7b5c6e3296196869066bf333/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/a5d079c2-1396-4f01-82f7-80769c44f640.d99f3028deb065a2bc07474fd2c3b3d2b380953e.de-de.xlf", "", "", "a5d079c2-1396-4f01-82f7-80769c44f640.d99f3028deb065a2bc07474fd2c3b3d2b380953e.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(8,1), "https://github.com/OpenLocalizationTest/oltest/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/e2e/f2c27ab1-8f74-4493-98ba-a7649837cbf6.md", "", "", "f2c27ab1-8f74-4493-98ba-a7649837cbf6.md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(8,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/f2c27ab1-8f74-4493-98ba-a7649837cbf6.a45b7595bce6941660e6dd85dc4c102f066ae408.de-de.xlf", "", "", "f2c27ab1-8f74-4493-98ba-a7649837cbf6.a45b7595bce6941660e6dd85dc4c102f066ae408.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Cells.Item(9,1), "https://github.com/OpenLocalizationTest/oltest/blob/a45b7595bce6941660e6dd85dc4c102f066ae408/.localization-config", "", "", ".localization-config")
